$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 849.0769
$ws.Range("I11").Value = 849.0769
$ws.Range("K11").Value = 849.0769
$ws.Range("M11").Value = -709.0769

$ws.Range("H69").Value = 2000
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = ""

$ws.Range("H72").Value = 2000
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = ""

$ws.Range("H88").Value = 1979.1765
$ws.Range("I88").Value = 2300
$ws.Range("J88").Value = 1936.4
$ws.Range("K88").Value = 2300
$ws.Range("L88").Value = 1936.4
$ws.Range("M88").Value = -1894
$ws.Range("N88").Value = -2748.4

$ws.Range("H91").Value = 1979.1765
$ws.Range("I91").Value = 2300
$ws.Range("J91").Value = 1936.4
$ws.Range("K91").Value = 2300
$ws.Range("L91").Value = 1936.4
$ws.Range("M91").Value = -896
$ws.Range("N91").Value = -4744.4

$ws.Range("H100").Value = 3802.5715
$ws.Range("I100").Value = 3956.6667
$ws.Range("J100").Value = 3687
$ws.Range("K100").Value = 3956.6667
$ws.Range("L100").Value = 3687
$ws.Range("M100").Value = -3415.6667
$ws.Range("N100").Value = -4769

$ws.Range("H125").Value = 5073.615
$ws.Range("I125").Value = 4876.3335
$ws.Range("J125").Value = 5132.8
$ws.Range("K125").Value = 43887.0015
$ws.Range("L125").Value = 46195.2
$ws.Range("M125").Value = -41427.0015
$ws.Range("N125").Value = -51115.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12823.017
$ws.Range("I32").Value = 14353.3
$ws.Range("K32").Value = 14353.3
$ws.Range("M32").Value = -14066.3

$ws.Range("H88").Value = 2425.8635
$ws.Range("I88").Value = 2185.5715
$ws.Range("J88").Value = 2538
$ws.Range("K88").Value = 2185.5715
$ws.Range("L88").Value = 2538
$ws.Range("M88").Value = -1779.5715
$ws.Range("N88").Value = -3350

$ws.Range("H91").Value = 2425.8635
$ws.Range("I91").Value = 2185.5715
$ws.Range("J91").Value = 2538
$ws.Range("K91").Value = 2185.5715
$ws.Range("L91").Value = 2538
$ws.Range("M91").Value = -781.5715
$ws.Range("N91").Value = -5346

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 205401.2
$ws.Range("I86").Value = 2499.5
$ws.Range("K86").Value = 2499.5
$ws.Range("M86").Value = -1376.5

$ws.Range("H89").Value = 205401.2
$ws.Range("I89").Value = 2499.5
$ws.Range("K89").Value = 12497.5
$ws.Range("M89").Value = -6881.5

$ws.Range("H105").Value = 5955205.5
$ws.Range("I105").Value = 9526463
$ws.Range("J105").Value = 3110.4443
$ws.Range("K105").Value = 9526463
$ws.Range("L105").Value = 3110.4443
$ws.Range("M105").Value = -9524716
$ws.Range("N105").Value = -6604.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 58388.332
$ws.Range("I62").Value = 65061.875
$ws.Range("K62").Value = 65061.875
$ws.Range("M62").Value = -64437.875

$ws.Range("H65").Value = 58388.332
$ws.Range("I65").Value = 65061.875
$ws.Range("K65").Value = 325309.375
$ws.Range("M65").Value = -322189.375

$ws.Range("H99").Value = 2548.75
$ws.Range("I99").Value = 2600
$ws.Range("J99").Value = 2497.5
$ws.Range("K99").Value = 2600
$ws.Range("L99").Value = 2497.5
$ws.Range("M99").Value = -1102
$ws.Range("N99").Value = -5493.5

$ws.Range("H105").Value = 1855.3846
$ws.Range("I105").Value = 1888.75
$ws.Range("J105").Value = 1802
$ws.Range("K105").Value = 1888.75
$ws.Range("L105").Value = 1802
$ws.Range("M105").Value = -141.75
$ws.Range("N105").Value = -5296

$ws.Range("H126").Value = 2548.75
$ws.Range("I126").Value = 2600
$ws.Range("J126").Value = 2497.5
$ws.Range("K126").Value = 7800
$ws.Range("L126").Value = 7492.5
$ws.Range("M126").Value = -5330
$ws.Range("N126").Value = -12432.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 6685.8335
$ws.Range("J35").Value = 6685.8335
$ws.Range("L35").Value = 20057.5005
$ws.Range("N35").Value = -20633.5005

$ws.Range("H36").Value = 6699.6665
$ws.Range("I36").Value = 100
$ws.Range("J36").Value = 9999.5
$ws.Range("K36").Value = 300
$ws.Range("L36").Value = 29998.5
$ws.Range("M36").Value = -131
$ws.Range("N36").Value = -30336.5

$ws.Range("H122").Value = 934.3333
$ws.Range("I122").Value = 304
$ws.Range("J122").Value = 1249.5
$ws.Range("K122").Value = 2736
$ws.Range("L122").Value = 11245.5
$ws.Range("M122").Value = -286
$ws.Range("N122").Value = -16145.5

$ws.Range("H125").Value = 3493.6365
$ws.Range("J125").Value = 3790
$ws.Range("L125").Value = 11370
$ws.Range("N125").Value = -21210

$ws.Range("H129").Value = 3572233.8
$ws.Range("I129").Value = 367.5
$ws.Range("J129").Value = 5000980
$ws.Range("K129").Value = 1102.5
$ws.Range("L129").Value = 15002940
$ws.Range("M129").Value = 3897.5
$ws.Range("N129").Value = -15012940

$ws.Range("H131").Value = 14547.328
$ws.Range("I131").Value = 344.05884
$ws.Range("J131").Value = 19376.44
$ws.Range("K131").Value = 1032.17652
$ws.Range("L131").Value = 58129.31999999999
$ws.Range("M131").Value = 4007.82348
$ws.Range("N131").Value = -68209.31999999999

$ws.Range("H133").Value = 4520
$ws.Range("J133").Value = 5971.5835
$ws.Range("L133").Value = 17914.7505
$ws.Range("N133").Value = -28034.7505

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3375
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 4500
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 4500
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -6496

$ws.Range("H83").Value = 3375
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 4500
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 22500
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -32484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6700.5557
$ws.Range("I7").Value = 4460
$ws.Range("K7").Value = 4460
$ws.Range("M7").Value = -4348

$ws.Range("H35").Value = 3671.8333
$ws.Range("I35").Value = 3671.8333
$ws.Range("K35").Value = 3671.8333
$ws.Range("M35").Value = -3335.8333

$ws.Range("H68").Value = 2583
$ws.Range("I68").Value = 2536.6667
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 2536.6667
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -1787.6667
$ws.Range("N68").Value = -4498

$ws.Range("H71").Value = 2583
$ws.Range("I71").Value = 2536.6667
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 12683.3335
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -8939.3335
$ws.Range("N71").Value = -22488

$ws.Range("H82").Value = 1755.2222
$ws.Range("I82").Value = 1460.4
$ws.Range("J82").Value = 2123.75
$ws.Range("K82").Value = 1460.4
$ws.Range("L82").Value = 2123.75
$ws.Range("M82").Value = -1099.4
$ws.Range("N82").Value = -2845.75

$ws.Range("H85").Value = 1755.2222
$ws.Range("I85").Value = 1460.4
$ws.Range("J85").Value = 2123.75
$ws.Range("K85").Value = 1460.4
$ws.Range("L85").Value = 2123.75
$ws.Range("M85").Value = -212.4000000000001
$ws.Range("N85").Value = -4619.75

$ws.Range("H126").Value = 6700.5557
$ws.Range("I126").Value = 4460
$ws.Range("K126").Value = 13380
$ws.Range("M126").Value = -10910

$ws.Range("H136").Value = 3207595.8
$ws.Range("I136").Value = 5436336.5
$ws.Range("J136").Value = 3780.9375
$ws.Range("K136").Value = 16309009.5
$ws.Range("L136").Value = 11342.8125
$ws.Range("M136").Value = -16306459.5
$ws.Range("N136").Value = -16442.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 25941.25
$ws.Range("I17").Value = 17921.666
$ws.Range("J17").Value = 50000
$ws.Range("K17").Value = 17921.666
$ws.Range("L17").Value = 50000
$ws.Range("M17").Value = -17749.666
$ws.Range("N17").Value = -50344

$ws.Range("H54").Value = 22643.334
$ws.Range("I54").Value = 14986.667
$ws.Range("J54").Value = 30300
$ws.Range("K54").Value = 14986.667
$ws.Range("L54").Value = 30300
$ws.Range("M54").Value = -14466.667
$ws.Range("N54").Value = -31340

$ws.Range("H57").Value = 46500
$ws.Range("I57").Value = 51000
$ws.Range("J57").Value = 42000
$ws.Range("K57").Value = 51000
$ws.Range("L57").Value = 42000
$ws.Range("M57").Value = -50246
$ws.Range("N57").Value = -43508

